# Add validation ("test") to Data.xlsx, auto_mlr
# - Duplicate the "Data" sheet into a new "Train Data" sheet (placed first)
# - Rename the original "Data" sheet to "Test Data"
# - Give both sheets the new set of column headers (Run, Roll temperature,
#   Target density, Calculated target porosity, C/5-6, 10C)
# - Trim "Train Data" down to 16 rows (header + 15 data rows) and
#   "Test Data" down to 4 rows (header + 3 data rows)
# - Make the "Train Data" comment on A1 its own (separate) threaded comment
# - Refresh view/selection state and make "Responses" the active tab

$wb = $excel.ActiveWorkbook

$dataWs  = $wb.Worksheets.Item("Data")
$firstWs = $wb.Worksheets.Item(1)

# Duplicate "Data" and put the copy in front of the first sheet -> becomes
# the new "Train Data" sheet (fresh sheetId, same comments/format as Data).
$dataWs.Copy($firstWs)

$trainWs = $wb.Worksheets.Item(1)
$trainWs.Name = "Train Data"

$testWs = $wb.Worksheets.Item("Data")
$testWs.Name = "Test Data"

$dpWs   = $wb.Worksheets.Item("Design Parameters")
$respWs = $wb.Worksheets.Item("Responses")

# New shared column headers for both the Train Data and Test Data sheets.
$headers = @("Run", "Roll temperature (°C)", "Target density (g/cm3)", "Calculated target porosity (%)", "C/5-6 (mAh/g)", "10C (mAh/g)")
$cols = @("A", "B", "C", "D", "E", "F")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $trainWs.Range($cols[$i] + "1").Value = $headers[$i]
    $testWs.Range($cols[$i] + "1").Value = $headers[$i]
}

# Header row grows to fit the wrapped column titles.
$trainWs.Rows(1).RowHeight = 58
$testWs.Rows(1).RowHeight = 58

# Trim row counts: Train Data keeps 15 data rows (2-16), Test Data keeps 3
# data rows (2-4) -- the remaining rows from the original 19-row sheet are
# removed.
$trainWs.Rows("17:19").Delete()
$testWs.Rows("5:19").Delete()

# The copied sheet shares the exact same threaded comment as the original
# "Data" sheet (now "Test Data"); give "Train Data" its own distinct
# comment instance (same text, new identity).
$trainComment = $trainWs.Range("A1").Comment
$commentText = $trainComment.Text()
$trainComment.Delete()
$trainWs.Range("A1").AddCommentThreaded($commentText)

# Refresh the selection / scroll state on each sheet.
$trainWs.Activate()
$trainWs.Rows("2:17").Select()

$testWs.Activate()
$testWs.Rows("2:4").Select()

$dpWs.Activate()
$dpWs.Rows("2:4").Select()

$respWs.Activate()
$respWs.Rows("2:4").Select()

Write-Host "Workbook restructured: Train Data / Test Data / Design Parameters / Responses"
